$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("Q:Q").Insert()
Write-Host ("UsedRange: " + $ws.UsedRange.Address())
Write-Host ("Q1: " + $ws.Range("Q1").Value)
Write-Host ("R1: " + $ws.Range("R1").Value)
Write-Host ("R61: " + $ws.Range("R61").Value)
Write-Host ("Q61 type: " + $ws.Range("Q61").Value)
